$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns keep their text formatting so that
# numeric-looking values (e.g. "0.999", "1.00") are stored as text,
# matching the original inlineStr cell content.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '56.713.29'
$ws.Range('E2').Value = '  +11.22%  '

# Row 3
$ws.Range('D3').Value = '3.254.65'
$ws.Range('E3').Value = '  +6.91%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').Value = '398.04'
$ws.Range('E5').Value = '  +2.53%  '

# Row 6
$ws.Range('D6').Value = '111.15'
$ws.Range('E6').Value = '  +10.66%  '

# Row 7
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').Value = '  +5.47%  '

# Row 8
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.09%  '

# Row 9
$ws.Range('D9').Value = '0.617'
$ws.Range('E9').Value = '  +6.96%  '

# Row 10
$ws.Range('D10').Value = '39.48'
$ws.Range('E10').Value = '  +8.02%  '

# Row 11
$ws.Range('D11').Value = '0.0949'
$ws.Range('E11').Value = '  +12.56%  '

# Row 12
$ws.Range('D12').Value = '0.142'
$ws.Range('E12').Value = '  +2.48%  '

# Row 13
$ws.Range('D13').Value = '3.767.54'
$ws.Range('E13').Value = '  +6.68%  '

# Row 14
$ws.Range('D14').Value = '19.25'
$ws.Range('E14').Value = '  +5.86%  '

# Row 15
$ws.Range('D15').Value = '8.10'
$ws.Range('E15').Value = '  +6.62%  '

# Row 16
$ws.Range('D16').Value = '3.247.34'
$ws.Range('E16').Value = '  +6.54%  '

# Row 17
$ws.Range('E17').Value = '  +5.51%  '

# Row 18
$ws.Range('D18').Value = '10.97'
$ws.Range('E18').Value = '  +3.35%  '

# Row 19
$ws.Range('D19').Value = '56.560.71'
$ws.Range('E19').Value = '  +10.79%  '

# Row 20
$ws.Range('E20').Value = '  +5.49%  '

# Row 21
$ws.Range('D21').Value = '0.0000104'
$ws.Range('E21').Value = '  +8.91%  '

# Row 22
$ws.Range('D22').Value = '13.01'
$ws.Range('E22').Value = '  +6.49%  '

# Row 23
$ws.Range('D23').Value = '302.65'
$ws.Range('E23').Value = '  +15.15%  '

# Row 24
$ws.Range('D24').Value = '74.94'
$ws.Range('E24').Value = '  +7.91%  '

# Row 25
$ws.Range('D25').Value = '3.20'
$ws.Range('E25').Value = '  +2.77%  '

# Row 26
$ws.Range('E26').Value = '  +3.37%  '

# Row 27
$ws.Range('D27').Value = '28.17'
$ws.Range('E27').Value = '  +6.07%  '

# Row 28
$ws.Range('E28').Value = '  +5.30%  '

# Row 29
$ws.Range('E29').Value = '  +2.71%  '

# Row 30
$ws.Range('D30').Value = '0.170'
$ws.Range('E30').Value = '  +6.01%  '

# Row 31
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.04%  '

# Row 32
$ws.Range('E32').Value = '  +6.75%  '

# Row 33
$ws.Range('D33').Value = '11.06'
$ws.Range('E33').Value = '  +5.79%  '

# Row 34
$ws.Range('D34').Value = '38.15'
$ws.Range('E34').Value = '  +6.77%  '

# Row 35
$ws.Range('D35').Value = '0.0486'
$ws.Range('E35').Value = '  +0.72%  '

# Row 36
$ws.Range('E36').Value = '  +3.50%  '

# Row 37
$ws.Range('D37').Value = '51.75'
$ws.Range('E37').Value = '  +3.57%  '

# Row 38
$ws.Range('E38').Value = '  +29.14%  '

# Row 39
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.14%  '

# Row 40
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '3.51'
$ws.Range('E40').Value = '  +6.30%  '

# Row 41
$ws.Range('D41').Value = '17.69'
$ws.Range('E41').Value = '  +7.48%  '

# Row 42
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = '134.72'
$ws.Range('E42').Value = '  +4.87%  '

# Row 43
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '1.94'
$ws.Range('E43').Value = '  +6.63%  '

# Row 44
$ws.Range('D44').Value = '0.119'
$ws.Range('E44').Value = '  +4.42%  '

# Row 45
$ws.Range('D45').Value = '3.96'
$ws.Range('E45').Value = '  +6.65%  '

# Row 46
$ws.Range('D46').Value = '0.284'
$ws.Range('E46').Value = '  -1.26%  '

# Row 47
$ws.Range('D47').Value = '22.13'
$ws.Range('E47').Value = '  +2.35%  '

# Row 48
$ws.Range('D48').Value = '2.146.37'
$ws.Range('E48').Value = '  +4.31%  '

# Row 49
$ws.Range('E49').Value = '  +1.62%  '

# Row 50
$ws.Range('D50').Value = '2.39'
$ws.Range('E50').Value = '  -3.33%  '

# Row 51
$ws.Range('D51').Value = '2.01'
$ws.Range('E51').Value = '  +42.00%  '
